$d = $word.ActiveDocument

function Set-ParaBold($paraIndex) {
    $p = $d.Paragraphs($paraIndex)
    $p.Range.Font.Bold = 1
}

function Set-ParaRunsXml($paraIndex, $xmlBody) {
    # Replaces the run content of the given paragraph (not including the
    # paragraph mark) with the supplied <w:r>/<w:proofErr> markup.
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r2 = $d.Range($r.Start, $r.End - 1)
    $xmlSnippet = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $xmlBody + '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part>' +
        '</pkg:package>'
    $r2.InsertXML($xmlSnippet)
}

# ---------------------------------------------------------------------
# 1. "Bug" heading becomes bold
# ---------------------------------------------------------------------
Set-ParaBold 1

# ---------------------------------------------------------------------
# 2. ".priority:[NONE, LOW, MEDIUM, HIGH]" split with a gramStart/gramEnd
#    proofErr bracketing ":["
# ---------------------------------------------------------------------
$body7 = (
    '<w:r><w:t>.priority</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>:[</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>NONE, LOW, MEDIUM, HIGH]</w:t></w:r>'
)
Set-ParaRunsXml 7 $body7

# ---------------------------------------------------------------------
# 3-5. ".dateCreated" / ".dateResolved" / ".dateClosed" split with a
#      spellStart/spellEnd proofErr bracketing the camelCase word
# ---------------------------------------------------------------------
$body8 = (
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>dateCreated</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)
Set-ParaRunsXml 8 $body8

$body9 = (
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>dateResolved</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)
Set-ParaRunsXml 9 $body9

$body10 = (
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>dateClosed</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)
Set-ParaRunsXml 10 $body10

# ---------------------------------------------------------------------
# 6. ".type: [BUG,ERROR]" split with a gramStart/gramEnd proofErr
#    bracketing ",ERROR"
# ---------------------------------------------------------------------
$body11 = (
    '<w:r><w:t>.type: [BUG</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>,ERROR</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>]</w:t></w:r>'
)
Set-ParaRunsXml 11 $body11

# ---------------------------------------------------------------------
# 7. ".status: [OPEN, INPROGRESS,RESOLVED,NEEDINFO,CLOSED]" split with a
#    gramStart/gramEnd proofErr bracketing ",RESOLVED,NEEDINFO,CLOSED"
# ---------------------------------------------------------------------
$body12 = (
    '<w:r><w:t>.status: [OPEN, INPROGRESS</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>,RESOLVED,NEEDINFO,CLOSED</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>]</w:t></w:r>'
)
Set-ParaRunsXml 12 $body12

# ---------------------------------------------------------------------
# 8. ".workHistory" split with a spellStart/spellEnd proofErr
# ---------------------------------------------------------------------
$body14 = (
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>workHistory</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)
Set-ParaRunsXml 14 $body14

# ---------------------------------------------------------------------
# 9. "User" heading becomes bold
# ---------------------------------------------------------------------
Set-ParaBold 17

# ---------------------------------------------------------------------
# 10. Insert a new ".loginName" paragraph right before the ".role" line
# ---------------------------------------------------------------------
$pRole = $d.Paragraphs(18)
$insertPoint = $d.Range($pRole.Range.Start, $pRole.Range.Start)
$insertPoint.InsertBefore(".loginName`r")

# The newly created paragraph is now paragraph 18; ".role" shifted to 19
$body18 = (
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>loginName</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)
Set-ParaRunsXml 18 $body18

# ---------------------------------------------------------------------
# 11. ".role: [ADMIN,PROJECTMANAGER,REPORTER,DEVELOPER]" split with a
#     gramStart/gramEnd proofErr bracketing ",PROJECTMANAGER,REPORTER,DEVELOPER"
# ---------------------------------------------------------------------
$body19 = (
    '<w:r><w:t>.role: [ADMIN</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>,PROJECTMANAGER,REPORTER,DEVELOPER</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>]</w:t></w:r>'
)
Set-ParaRunsXml 19 $body19

# ---------------------------------------------------------------------
# 12. "Project" heading becomes bold
# ---------------------------------------------------------------------
Set-ParaBold 21

Write-Host "Edits applied successfully"
